$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 11985.571
$ws.Range("I21").Value = 1950
$ws.Range("J21").Value = 15999.8
$ws.Range("K21").Value = 1950
$ws.Range("L21").Value = 15999.8
$ws.Range("M21").Value = -1482
$ws.Range("N21").Value = -16935.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 11985.571
$ws.Range("I23").Value = 1950
$ws.Range("J23").Value = 15999.8
$ws.Range("K23").Value = 1950
$ws.Range("L23").Value = 15999.8
$ws.Range("M23").Value = -1716
$ws.Range("N23").Value = -16467.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 108.28
$ws.Range("I33").Value = 117.09091
$ws.Range("J33").Value = 101.35714
$ws.Range("K33").Value = 117.09091
$ws.Range("L33").Value = 101.35714
$ws.Range("M33").Value = 111.90909
$ws.Range("N33").Value = -559.35714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 245541.33
$ws.Range("J95").Value = 245541.33
$ws.Range("L95").Value = 245541.33
$ws.Range("N95").Value = -251033.33

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 584067.5600000001
$ws.Range("J132").Value = 28751.25
$ws.Range("L132").Value = 86253.75
$ws.Range("N132").Value = -91313.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 35715380
$ws.Range("I137").Value = 55556624
$ws.Range("J137").Value = 1140.6
$ws.Range("K137").Value = 166669872
$ws.Range("L137").Value = 3421.8
$ws.Range("M137").Value = -166667322
$ws.Range("N137").Value = -8521.799999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4151.5
$ws.Range("I141").Value = 2914.3333
$ws.Range("J141").Value = 7863
$ws.Range("K141").Value = 8742.999899999999
$ws.Range("L141").Value = 23589
$ws.Range("M141").Value = -3562.999899999999
$ws.Range("N141").Value = -33949

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 4000
$ws.Range("J17").Value = 4000
$ws.Range("L17").Value = 4000
$ws.Range("N17").Value = -4346

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H64").Value = 25000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 25000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 25000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -25496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H67").Value = 25000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 25000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 25000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -26716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3046.639
$ws.Range("I134").Value = 1722.826
$ws.Range("J134").Value = 5388.769
$ws.Range("K134").Value = 5168.478
$ws.Range("L134").Value = 16166.307
$ws.Range("M134").Value = -2633.478
$ws.Range("N134").Value = -21236.307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1008.2222
$ws.Range("I16").Value = 1020.125
$ws.Range("J16").Value = 913
$ws.Range("K16").Value = 1020.125
$ws.Range("L16").Value = 913
$ws.Range("M16").Value = -733.125
$ws.Range("N16").Value = -1487

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4440.079
$ws.Range("I31").Value = 1157.4667
$ws.Range("J31").Value = 16749.875
$ws.Range("K31").Value = 1157.4667
$ws.Range("L31").Value = 16749.875
$ws.Range("M31").Value = -862.4666999999999
$ws.Range("N31").Value = -17339.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4440.079
$ws.Range("I34").Value = 1157.4667
$ws.Range("J34").Value = 16749.875
$ws.Range("K34").Value = 1157.4667
$ws.Range("L34").Value = 16749.875
$ws.Range("M34").Value = -955.4666999999999
$ws.Range("N34").Value = -17153.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1008.2222
$ws.Range("I113").Value = 1020.125
$ws.Range("J113").Value = 913
$ws.Range("K113").Value = 1020.125
$ws.Range("L113").Value = 913
$ws.Range("M113").Value = 1149.875
$ws.Range("N113").Value = -5253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2334.2222
$ws.Range("I132").Value = 2005
$ws.Range("K132").Value = 6015
$ws.Range("M132").Value = -3485

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2456.611
$ws.Range("I134").Value = 1380.2693
$ws.Range("K134").Value = 4140.8079
$ws.Range("M134").Value = -1605.8079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 825.2857
$ws.Range("I97").Value = 1098.25
$ws.Range("J97").Value = 461.33334
$ws.Range("K97").Value = 3294.75
$ws.Range("L97").Value = 1384.00002
$ws.Range("M97").Value = -2798.75
$ws.Range("N97").Value = -2376.00002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1221.8572
$ws.Range("I129").Value = 524
$ws.Range("J129").Value = 2966.5
$ws.Range("K129").Value = 1572
$ws.Range("L129").Value = 8899.5
$ws.Range("M129").Value = 3428
$ws.Range("N129").Value = -18899.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5129491.5
$ws.Range("I131").Value = 518.8
$ws.Range("J131").Value = 5556905.5
$ws.Range("K131").Value = 1556.4
$ws.Range("L131").Value = 16670716.5
$ws.Range("M131").Value = 3483.6
$ws.Range("N131").Value = -16680796.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3583.182
$ws.Range("I139").Value = 3439.0908
$ws.Range("J139").Value = 3727.2727
$ws.Range("K139").Value = 10317.2724
$ws.Range("L139").Value = 11181.8181
$ws.Range("M139").Value = -5177.2724
$ws.Range("N139").Value = -21461.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 6284.4814
$ws.Range("I140").Value = 14139.875
$ws.Range("J140").Value = 2976.9473
$ws.Range("K140").Value = 42419.625
$ws.Range("L140").Value = 8930.841899999999
$ws.Range("M140").Value = -37239.625
$ws.Range("N140").Value = -19290.8419

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 24615.166
$ws.Range("J57").Value = 23540.2
$ws.Range("L57").Value = 23540.2
$ws.Range("N57").Value = -25180.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 11000
$ws.Range("J92").Value = 11000
$ws.Range("L92").Value = 11000
$ws.Range("N92").Value = -14744

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3289.484
$ws.Range("I132").Value = 2600.5789
$ws.Range("K132").Value = 7801.736699999999
$ws.Range("M132").Value = -5271.736699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3997.9412
$ws.Range("I132").Value = 3022.8333
$ws.Range("K132").Value = 9068.499899999999
$ws.Range("M132").Value = -6538.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2035.1224
$ws.Range("I132").Value = 1758.4828
$ws.Range("J132").Value = 2436.25
$ws.Range("K132").Value = 5275.4484
$ws.Range("L132").Value = 7308.75
$ws.Range("M132").Value = -2745.4484
$ws.Range("N132").Value = -12368.75
